# Apply "Added advise on PROGn, INITn, DONE pins" edit.
$wb = $excel.ActiveWorkbook
$wsSchematic = $wb.Worksheets.Item("Schematic")
$wsLayout = $wb.Worksheets.Item("Layout")

# --- Schematic sheet ---

# Expand the DONE-pin note to also call out INITn / PROGn.
$wsSchematic.Range("B19").Value = "INITn, PROGn, DONE pin should not be pulled low during FPGA boot up, this is an open drain."

# Update the SoM current-draw note.
$wsSchematic.Range("B4").Value = "1.8V/3.3V sourced from SoM consumes no more than 800mA"

# Fix "VBUS conencted to USB 5V" typo -> "VBUS connected to USB 5V".
$wsSchematic.Range("B3").Value = "VBUS connected to USB 5V"

# Insert a new row for the "Power enable" review item right after the
# "5V USB input hooked up" row, pushing everything below it down by one.
[void]$wsSchematic.Rows.Item(3).Insert()
$wsSchematic.Range("B3").Value = "Power enable must be high to turn on the SoM"

# "Bank voltage allocation" section header (now row 6) becomes bold.
$wsSchematic.Range("A6").Font.Bold = $true

# "MIPI" section header (now row 24) becomes bold.
$wsSchematic.Range("A24").Font.Bold = $true

# Leave the cursor where the edits were last made.
[void]$wsSchematic.Range("A24").Select()

# --- Layout sheet ---

# "Check MIPI intra pair length matches" -> "...match"
$wsLayout.Range("B2").Value = "Check MIPI intra pair length match"

# Layout is the sheet left active/selected when the workbook was saved.
[void]$wsLayout.Select()
[void]$wsLayout.Range("B3").Select()
